$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 30
    3  = 30
    4  = 30
    5  = 31
    6  = 31
    7  = 31
    8  = 31
    9  = 31
    10 = 31
    11 = 31
    12 = 32
    13 = 32
    14 = 32
    15 = 33
    16 = 33
    17 = 34
    18 = 34
    19 = 34
    20 = 34
    21 = 34
    22 = 34
    23 = 34
    24 = 34
    25 = 34
    26 = 34
    27 = 34
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 2).Value = $values[$row]
}
